# Extending chb06 interictal test set
# (commit message says chb06, but the actual data change is new chb10
#  Test/Interictal rows — following the authoritative xml diff.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Turn the flat G42/G43 values into a "=900*2" formula (matches
#        the pattern already used by rows 44/45 just below them). ---
$ws.Range("G42").Formula = "=900*2"
$ws.Range("G43").Formula = "=900*2"

# --- 2. Append 9 new interictal test rows for chb10 (rows 52-60). ---
$newFiles = @(
    "chb10_07.edf",
    "chb10_08.edf",
    "chb10_13.edf",
    "chb10_14.edf",
    "chb10_15.edf",
    "chb10_16.edf",
    "chb10_17.edf",
    "chb10_18.edf",
    "chb10_19.edf"
)

$startRow = 52
for ($i = 0; $i -lt $newFiles.Count; $i++) {
    $r = $startRow + $i

    # Copy the formatting of the row directly above (row 51) down onto
    # column B so the new filename cells keep the same style (s="2") as
    # every other "Filename" cell in the table.
    $ws.Range("B51").Copy()
    $ws.Range("B" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value = "chb10"
    $ws.Range("B" + $r).Value = $newFiles[$i]
    $ws.Range("C" + $r).Value = "Test"
    $ws.Range("D" + $r).Value = "Interictal"
    $ws.Range("E" + $r).Value = 0
    $ws.Range("F" + $r).Value = 0
    $ws.Range("G" + $r).Value = 900
}

$excel.CutCopyMode = 0

# --- 3. Grow the hidden _FilterDatabase defined name to cover the new rows. ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$G`$60"
    }
}

# --- 4. Update the view: scroll down a row and move the active selection. ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H43").Select()
